$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PART")
$ws.Range("Z35").Value = "test"
$ws.Range("Z35").ClearFormats()
$ws.Range("Z35").Font.Name = "Monaco"
$ws.Range("Z35").Font.Size = 11
$ws.Range("Z35").Font.Color = 255
Write-Host "done"
